$wb = $excel.ActiveWorkbook

# Updated "想去人数" (number of people interested) counts need to be applied to
# both the "展览" sheet and the "全部类型" sheet, which contain duplicate data.
$updates = @{
    2 = 2194
    3 = 1664
    4 = 324
    5 = 1069
    6 = 676
    8 = 5760
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
